$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-07-09 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-10 Thursday", 2) | Out-Null

# Update each arithmetic-problem cell in the table, in row-major (document) order
$tbl = $d.Tables.Item(1)
$newValues = @(
    "96-63=33", "39-20=19", "43+21=64", "41+27=68", "30-24=6", "98-66=32", "71-61=10", "83-16=67",
    "17+15=32", "96-62=34", "7+8=15", "37+60=97", "76-28=48", "79+7=86", "98-30=68", "82+13=95",
    "66-16=50", "23+43=66", "33+59=92", "32+20=52", "42+6=48", "7+69=76", "91-88=3", "6+38=44",
    "54-43=11", "19+8=27", "7+68=75", "20+56=76", "38+40=78", "37-2=35", "20-19=1", "1+85=86",
    "20+31=51", "22+22=44", "14+8=22", "50-8=42", "60-28=32", "39+55=94", "74-74=0", "55-33=22",
    "56+22=78", "62-52=10", "60-40=20", "58+12=70", "4+29=33", "4+21=25", "73+2=75", "98+0=98",
    "35+8=43", "33+61=94", "63-43=20", "67-46=21", "2+84=86", "2+24=26", "49+11=60", "95+1=96",
    "13+1=14", "23+21=44", "31+4=35", "95-77=18", "37+42=79", "78-9=69", "50+49=99", "39-7=32",
    "99-82=17", "41+14=55", "77-59=18", "1+39=40", "83-47=36", "94-65=29", "26-19=7", "85-73=12",
    "54-42=12", "39+41=80", "11+52=63", "88-40=48", "28+46=74", "3+40=43", "1+83=84", "67-3=64",
    "13+83=96", "57+5=62", "86-66=20", "21+6=27", "63+32=95", "79-35=44", "85-57=28", "28+35=63",
    "26+18=44", "64+11=75", "49-42=7", "17-12=5", "46-36=10", "27+52=79", "20+23=43", "68-59=9",
    "51+23=74", "36+19=55", "20-10=10", "19+58=77"
)

$rows = $tbl.Rows.Count
$cols = $tbl.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $tbl.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Output "Updated $idx cells"
